$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = " 16:43"
$ws.Range("D2").Value = " 24.5 °C"
$ws.Range("F2").Value = " 41 %"

# Row 3
$ws.Range("B3").Value = " 16:43"
$ws.Range("D3").Value = " 26 °C"
$ws.Range("E3").Value = " 29 °C"
$ws.Range("F3").Value = " 41 %"

# Row 4
$ws.Range("B4").Value = " 16:43"
$ws.Range("D4").Value = " 24.2 °C"
$ws.Range("E4").Value = " 27 °C"
$ws.Range("F4").Value = " 46 %"

# Row 5
$ws.Range("B5").Value = " 16:43"
$ws.Range("D5").Value = " 25 °C"
$ws.Range("F5").Value = " 40 %"

# Row 6
$ws.Range("B6").Value = " 16:43"
$ws.Range("F6").Value = " 41 %"

# Row 7
$ws.Range("B7").Value = " 16:43"
$ws.Range("D7").Value = " 25 °C"
$ws.Range("E7").Value = " 25 °C"
$ws.Range("F7").Value = " 37 %"

# Row 8
$ws.Range("B8").Value = " 16:43"
$ws.Range("D8").Value = " 23 °C"
$ws.Range("F8").Value = " 49 %"

# Row 9
$ws.Range("B9").Value = " 16:43"
$ws.Range("E9").Value = " 26 °C"
$ws.Range("F9").Value = " 40 %"

# Row 10
$ws.Range("B10").Value = " 16:43"
$ws.Range("F10").Value = " 40 %"

# Row 11
$ws.Range("B11").Value = " 16:43"
$ws.Range("D11").Value = " 25.6 °C"
$ws.Range("E11").Value = " 27 °C"
$ws.Range("F11").Value = " 36 %"

# Row 12
$ws.Range("B12").Value = " 16:43"
$ws.Range("D12").Value = " 24.4 °C"
$ws.Range("F12").Value = " 39 %"

# Row 13
$ws.Range("B13").Value = " 16:43"
$ws.Range("D13").Value = " 25.6 °C"
$ws.Range("E13").Value = " 27 °C"
$ws.Range("F13").Value = " 36 %"

# Row 14
$ws.Range("B14").Value = " 16:43"
$ws.Range("D14").Value = " 22.8 °C"

# Row 15
$ws.Range("B15").Value = " 16:43"
$ws.Range("D15").Value = " 25.8 °C"
$ws.Range("E15").Value = " 28 °C"
$ws.Range("F15").Value = " 41 %"

# Row 16
$ws.Range("B16").Value = " 16:43"
$ws.Range("D16").Value = " 19 °C"
$ws.Range("F16").Value = " 42 %"

# Row 17
$ws.Range("B17").Value = " 16:43"
$ws.Range("D17").Value = " 25.6 °C"
$ws.Range("E17").Value = " 27 °C"
$ws.Range("F17").Value = " 36 %"

# Row 18
$ws.Range("B18").Value = " 16:43"
$ws.Range("D18").Value = " 23 °C"
$ws.Range("E18").Value = " 24 °C"
$ws.Range("F18").Value = " 52 %"

# Row 19
$ws.Range("B19").Value = " 16:43"
$ws.Range("D19").Value = " 26 °C"
$ws.Range("E19").Value = " 28 °C"
$ws.Range("F19").Value = " 41 %"

# Row 20
$ws.Range("B20").Value = " 16:43"
$ws.Range("D20").Value = " 26 °C"
$ws.Range("E20").Value = " 29 °C"
$ws.Range("F20").Value = " 41 %"
